$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header row values for new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style from an existing header cell (AC1) to the new headers
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Fill in Wins/Losses/Ties values for each data row (2-44)
for ($row = 2; $row -le 44; $row++) {
    $ws.Cells.Item($row, 30).Value = 71   # AD = column 30
    $ws.Cells.Item($row, 31).Value = 91   # AE = column 31
    $ws.Cells.Item($row, 32).Value = 0    # AF = column 32
}
